# Insert a new row of swim data (50 FR, 41.19r) at row 3, pushing existing
# rows 3-36 down to 4-37, matching meets recorded through December 2019.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 3 (shifts rows 3.. down by one)
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the new swim result
$ws.Range("A3").Value = "50 FR"
$ws.Range("B3").Value = "41.19r"
$ws.Range("C3").Value = 36557.50047673611
$ws.Range("C3").NumberFormat = "ss.00"
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 162
$ws.Range("F3").Value = """Slower than B"""
$ws.Range("G3").Value = "2019 CT RYWC Candlewood Cup"
$ws.Range("H3").Value = "CT"
$ws.Range("I3").Value = "PAC"
$ws.Range("J3").Value = "11/9/2019"
